$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the green header fill style to E1 (same fill as the other header cells)
$ws.Range("E1").Interior.Color = $ws.Range("B1").Interior.Color

# Move "has one coordinator thru staff courses supervision" from E8 to D6
$ws.Range("D6").Value = $ws.Range("E8").Value2
$ws.Range("D6").WrapText = $true

# Clear the old E8 cell (content + formatting) entirely
$ws.Range("E8").ClearContents()
$ws.Range("E8").ClearFormats()

# Row 6 now needs the 30pt height that row 8 used to have; row 8 goes back to default
$ws.Rows("6:6").RowHeight = 30
$ws.Rows("8:8").EntireRow.AutoFit()

# Update the active selection/cursor position
$ws.Range("F1").Select()
